$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date formatting (style) from an existing date cell so the new
# cells reuse the same cellXf (numFmtId 14) rather than creating a new one.
$ws.Range("A51").Copy()
$ws.Range("A52:A53").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 52: fill in the previously-empty cells and fix the formula
$ws.Range("A52").Value = 43742
$ws.Range("B52").Value = 2275.4950372970102
$ws.Range("D52").Formula = "=100*(B52-C52)/C52"
$ws.Range("E52").Value = 169

# Row 53: new CRM test entry from 10/04/2019
$ws.Range("A53").Value = 43742
$ws.Range("B53").Value = 2268.8014966576202
$ws.Range("C53").Value = 2207.0300000000002
$ws.Range("D53").Formula = "=100*(B53-C53)/C53"
$ws.Range("E53").Value = 169
$ws.Range("F53").Value = "opened crm (8/7/2019)"

$ws.Range("A53").Select()
